$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4662.88
$ws.Range("L51").Value = 4469.25
$ws.Range("N51").Value = -5437.25
$ws.Range("J51").Value = 4469.25
$ws.Range("K86").Value = 11514.857
$ws.Range("N86").Value = -83024.28999999999
$ws.Range("L86").Value = 80778.28999999999
$ws.Range("M86").Value = -10391.857
$ws.Range("I86").Value = 11514.857
$ws.Range("J86").Value = 80778.28999999999
$ws.Range("H86").Value = 34602.668
$ws.Range("N89").Value = -415123.45
$ws.Range("K89").Value = 57574.285
$ws.Range("M89").Value = -51958.285
$ws.Range("L89").Value = 403891.45
$ws.Range("H89").Value = 34602.668
$ws.Range("J89").Value = 80778.28999999999
$ws.Range("I89").Value = 11514.857
$ws.Range("J97").Value = 2057.2307
$ws.Range("N97").Value = -7163.6921
$ws.Range("L97").Value = 6171.6921
$ws.Range("H97").Value = 2057.2307
$ws.Range("J106").Value = 12000
$ws.Range("H106").Value = 4999.8
$ws.Range("L106").Value = 12000
$ws.Range("N106").Value = -13262
$ws.Range("M111").Value = -52.14259999999967
$ws.Range("K111").Value = 3119.1426
$ws.Range("H111").Value = 1159.75
$ws.Range("I111").Value = 1039.7142
$ws.Range("L112").Value = 177632.835
$ws.Range("J112").Value = 59210.945
$ws.Range("H112").Value = 94324.59
$ws.Range("N112").Value = -179848.835
$ws.Range("K132").Value = 7318.7724
$ws.Range("H132").Value = 2439.5908
$ws.Range("I132").Value = 2439.5908
$ws.Range("M132").Value = -4788.7724
$ws.Range("L137").Value = 9444.599999999999
$ws.Range("K137").Value = 6407.3181
$ws.Range("I137").Value = 2135.7727
$ws.Range("M137").Value = -3857.3181
$ws.Range("N137").Value = -14544.6
$ws.Range("J137").Value = 3148.2
$ws.Range("H137").Value = 2452.1562
$ws.Range("K141").Value = 2368.88892
$ws.Range("I141").Value = 789.62964
$ws.Range("H141").Value = 740.5
$ws.Range("M141").Value = 2811.11108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1839641.2
$ws.Range("M2").Value = -1839528.2
$ws.Range("K2").Value = 1839641.2
$ws.Range("H2").Value = 1338839
$ws.Range("M32").Value = -2074.2778
$ws.Range("I32").Value = 2361.2778
$ws.Range("H32").Value = 2467.8086
$ws.Range("K32").Value = 2361.2778
$ws.Range("M116").Value = -1837347.2
$ws.Range("I116").Value = 1839641.2
$ws.Range("H116").Value = 1338839
$ws.Range("K116").Value = 1839641.2
$ws.Range("L122").Value = 8321.25
$ws.Range("M122").Value = -5853.499899999999
$ws.Range("H122").Value = 2769.3125
$ws.Range("N122").Value = -13221.25
$ws.Range("I122").Value = 2767.8333
$ws.Range("J122").Value = 2773.75
$ws.Range("K122").Value = 8303.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -1839527.2
$ws.Range("K3").Value = 1839641.2
$ws.Range("I3").Value = 1839641.2
$ws.Range("H3").Value = 1338839
$ws.Range("M107").Value = 24.14280000000008
$ws.Range("K107").Value = 1895.8572
$ws.Range("H107").Value = 54559.95
$ws.Range("I107").Value = 1895.8572
$ws.Range("K113").Value = 8500
$ws.Range("M113").Value = -6330
$ws.Range("H113").Value = 8500
$ws.Range("I113").Value = 8500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 4999.6665
$ws.Range("L16").Value = 4999.6665
$ws.Range("M16").Value = -3623601.2
$ws.Range("N16").Value = -5573.6665
$ws.Range("H16").Value = 1814444
$ws.Range("I16").Value = 3623888.2
$ws.Range("K16").Value = 3623888.2
$ws.Range("H31").Value = 8509.823
$ws.Range("N31").Value = -11443.059
$ws.Range("L31").Value = 10853.059
$ws.Range("M31").Value = -5871.5884
$ws.Range("I31").Value = 6166.5884
$ws.Range("K31").Value = 6166.5884
$ws.Range("J31").Value = 10853.059
$ws.Range("M34").Value = -5964.5884
$ws.Range("I34").Value = 6166.5884
$ws.Range("N34").Value = -11257.059
$ws.Range("J34").Value = 10853.059
$ws.Range("H34").Value = 8509.823
$ws.Range("L34").Value = 10853.059
$ws.Range("K34").Value = 6166.5884
$ws.Range("H51").Value = 18499.875
$ws.Range("L51").Value = 18285.715
$ws.Range("N51").Value = -19757.715
$ws.Range("J51").Value = 18285.715
$ws.Range("H61").Value = 18499.875
$ws.Range("J61").Value = 18285.715
$ws.Range("L61").Value = 18285.715
$ws.Range("N61").Value = -18981.715
$ws.Range("K113").Value = 3623888.2
$ws.Range("M113").Value = -3621718.2
$ws.Range("H113").Value = 1814444
$ws.Range("I113").Value = 3623888.2
$ws.Range("L113").Value = 4999.6665
$ws.Range("N113").Value = -9339.666499999999
$ws.Range("J113").Value = 4999.6665
$ws.Range("K134").Value = 19235347.5
$ws.Range("H134").Value = 5209855
$ws.Range("L134").Value = 4506.6666
$ws.Range("J134").Value = 1502.2222
$ws.Range("M134").Value = -19232812.5
$ws.Range("I134").Value = 6411782.5
$ws.Range("N134").Value = -9576.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N2").Value = -7292.666800000001
$ws.Range("I2").Value = 92
$ws.Range("J2").Value = 1177.7778
$ws.Range("M2").Value = -439
$ws.Range("L2").Value = 7066.666800000001
$ws.Range("K2").Value = 552
$ws.Range("H2").Value = 790
$ws.Range("H107").Value = 841.16
$ws.Range("L107").Value = 2742.8334
$ws.Range("N107").Value = -6582.8334
$ws.Range("J107").Value = 914.2778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 8468.875
$ws.Range("L10").Value = 8999
$ws.Range("M10").Value = -8224.143
$ws.Range("K10").Value = 8393.143
$ws.Range("J10").Value = 8999
$ws.Range("N10").Value = -9337
$ws.Range("I10").Value = 8393.143
$ws.Range("J18").Value = 47900
$ws.Range("I18").Value = 48299.668
$ws.Range("H18").Value = 48199.75
$ws.Range("K18").Value = 48299.668
$ws.Range("L18").Value = 47900
$ws.Range("N18").Value = -48486
$ws.Range("M18").Value = -48006.668
$ws.Range("J100").Value = 129999
$ws.Range("N100").Value = -132163
$ws.Range("L100").Value = 129999
$ws.Range("H100").Value = 129999
$ws.Range("M102").Value = -105.8571999999999
$ws.Range("I102").Value = 1727.8572
$ws.Range("H102").Value = 7136.875
$ws.Range("J102").Value = 45000
$ws.Range("N102").Value = -48244
$ws.Range("K102").Value = 1727.8572
$ws.Range("L102").Value = 45000
$ws.Range("M122").Value = -4964.6362
$ws.Range("H122").Value = 3915.3257
$ws.Range("I122").Value = 2471.5454
$ws.Range("K122").Value = 7414.6362
$ws.Range("H126").Value = 6957.067
$ws.Range("M126").Value = -18886.1432
$ws.Range("I126").Value = 7118.7144
$ws.Range("K126").Value = 21356.1432
$ws.Range("K132").Value = 75004662
$ws.Range("H132").Value = 17858596
$ws.Range("I132").Value = 25001554
$ws.Range("M132").Value = -75002132

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J16").Value = 3573.9443
$ws.Range("L16").Value = 3573.9443
$ws.Range("M16").Value = -1197.1428
$ws.Range("N16").Value = -3913.9443
$ws.Range("H16").Value = 2608.4688
$ws.Range("I16").Value = 1367.1428
$ws.Range("K16").Value = 1367.1428
$ws.Range("K132").Value = 68724762
$ws.Range("H132").Value = 20887380
$ws.Range("I132").Value = 22908254
$ws.Range("M132").Value = -68722232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M14").Value = 167
$ws.Range("K14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("N96").Value = -4412
$ws.Range("J96").Value = 1666
$ws.Range("L96").Value = 1666
$ws.Range("H96").Value = 1849.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1776.6667
$ws.Range("H122").Value = 1568
$ws.Range("N122").Value = -13900
$ws.Range("I122").Value = 1408.8889
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4226.6667
$ws.Range("H136").Value = 16669015
$ws.Range("M136").Value = -51728772
$ws.Range("I136").Value = 17243774
$ws.Range("K136").Value = 51731322

Write-Host "Applied all updates"